# Auto-generated Excel COM-interop script
# Applies numeric value updates to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
# per the commit's profit-recalculation diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Cells.Item(32, 8).Value = 22670  # H32: was 20744.445
$ws.Cells.Item(32, 9).Value = 26666.666  # I32: was 20000
$ws.Cells.Item(32, 11).Value = 26666.666  # K32: was 20000
$ws.Cells.Item(32, 13).Value = -26340.666  # M32: was -19674
# Row 49
$ws.Cells.Item(49, 8).Value = 1246.3334  # H49: was 227
$ws.Cells.Item(49, 10).Value = 6316  # J49: was 200
$ws.Cells.Item(49, 12).Value = 18948  # L49: was 600
$ws.Cells.Item(49, 14).Value = -19220  # N49: was -872
# Row 51
$ws.Cells.Item(51, 8).Value = 2908.4167  # H51: was 2869.5652
$ws.Cells.Item(51, 9).Value = 2625  # I51: was 2531.25
$ws.Cells.Item(51, 10).Value = 3475.25  # J51: was 3642.8572
$ws.Cells.Item(51, 11).Value = 2625  # K51: was 2531.25
$ws.Cells.Item(51, 12).Value = 3475.25  # L51: was 3642.8572
$ws.Cells.Item(51, 13).Value = -2141  # M51: was -2047.25
$ws.Cells.Item(51, 14).Value = -4443.25  # N51: was -4610.8572
# Row 76
$ws.Cells.Item(76, 8).Value = 7026.3335  # H76: was 10000.5
$ws.Cells.Item(76, 9).Value = 9996.5  # I76: was 10000.5
$ws.Cells.Item(76, 10).Value = 5541.25  # J76: was 0
$ws.Cells.Item(76, 11).Value = 9996.5  # K76: was 10000.5
$ws.Cells.Item(76, 12).Value = 5541.25  # L76: was 0
$ws.Cells.Item(76, 13).Value = -9681.5  # M76: was -9685.5
$ws.Cells.Item(76, 14).Value = -6171.25  # N76: was None
# Row 79
$ws.Cells.Item(79, 8).Value = 7026.3335  # H79: was 10000.5
$ws.Cells.Item(79, 9).Value = 9996.5  # I79: was 10000.5
$ws.Cells.Item(79, 10).Value = 5541.25  # J79: was 0
$ws.Cells.Item(79, 11).Value = 9996.5  # K79: was 10000.5
$ws.Cells.Item(79, 12).Value = 5541.25  # L79: was 0
$ws.Cells.Item(79, 13).Value = -8904.5  # M79: was -8908.5
$ws.Cells.Item(79, 14).Value = -7725.25  # N79: was None
# Row 96
$ws.Cells.Item(96, 8).Value = 596.94116  # H96: was 523.5
$ws.Cells.Item(96, 9).Value = 619.3125  # I96: was 540.2353000000001
$ws.Cells.Item(96, 11).Value = 1857.9375  # K96: was 1620.7059
$ws.Cells.Item(96, 13).Value = -484.9375  # M96: was -247.7059000000002
# Row 97
$ws.Cells.Item(97, 8).Value = 166666  # H97: was 83887.5
$ws.Cells.Item(97, 10).Value = 166666  # J97: was 83887.5
$ws.Cells.Item(97, 12).Value = 499998  # L97: was 251662.5
$ws.Cells.Item(97, 14).Value = -500990  # N97: was -252654.5
# Row 99
$ws.Cells.Item(99, 8).Value = 5542.6  # H99: was 5609.8
$ws.Cells.Item(99, 9).Value = 821  # I99: was 899.5
$ws.Cells.Item(99, 10).Value = 12625  # J99: was 8750
$ws.Cells.Item(99, 11).Value = 2463  # K99: was 2698.5
$ws.Cells.Item(99, 12).Value = 37875  # L99: was 26250
$ws.Cells.Item(99, 13).Value = -965  # M99: was -1200.5
$ws.Cells.Item(99, 14).Value = -40871  # N99: was -29246
# Row 116
$ws.Cells.Item(116, 8).Value = 9455.27  # H116: was 9563.48
$ws.Cells.Item(116, 9).Value = 10819.667  # I116: was 11189.637
$ws.Cells.Item(116, 11).Value = 10819.667  # K116: was 11189.637
$ws.Cells.Item(116, 13).Value = -7377.666999999999  # M116: was -7747.637000000001
# Row 118
$ws.Cells.Item(118, 8).Value = 77895.46000000001  # H118: was 91982.91
$ws.Cells.Item(118, 9).Value = 77895.46000000001  # I118: was 91982.91
$ws.Cells.Item(118, 11).Value = 233686.38  # K118: was 275948.73
$ws.Cells.Item(118, 13).Value = -232029.38  # M118: was -274291.73
# Row 132
$ws.Cells.Item(132, 8).Value = 1295.6061  # H132: was 1363.7
$ws.Cells.Item(132, 9).Value = 1307.3572  # I132: was 1361.9615
$ws.Cells.Item(132, 10).Value = 1229.8  # J132: was 1375
$ws.Cells.Item(132, 11).Value = 3922.0716  # K132: was 4085.8845
$ws.Cells.Item(132, 12).Value = 3689.4  # L132: was 4125
$ws.Cells.Item(132, 13).Value = -1392.0716  # M132: was -1555.8845
$ws.Cells.Item(132, 14).Value = -8749.4  # N132: was -9185

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Cells.Item(74, 8).Value = 7960.654  # H74: was 9366.137000000001
$ws.Cells.Item(74, 9).Value = 961.0952  # I74: was 1085.0625
$ws.Cells.Item(74, 10).Value = 37358.8  # J74: was 31449
$ws.Cells.Item(74, 11).Value = 961.0952  # K74: was 1085.0625
$ws.Cells.Item(74, 12).Value = 37358.8  # L74: was 31449
$ws.Cells.Item(74, 13).Value = -87.09519999999998  # M74: was -211.0625
$ws.Cells.Item(74, 14).Value = -39106.8  # N74: was -33197
# Row 77
$ws.Cells.Item(77, 8).Value = 7960.654  # H77: was 9366.137000000001
$ws.Cells.Item(77, 9).Value = 961.0952  # I77: was 1085.0625
$ws.Cells.Item(77, 10).Value = 37358.8  # J77: was 31449
$ws.Cells.Item(77, 11).Value = 4805.476  # K77: was 5425.3125
$ws.Cells.Item(77, 12).Value = 186794  # L77: was 157245
$ws.Cells.Item(77, 13).Value = -437.4759999999997  # M77: was -1057.3125
$ws.Cells.Item(77, 14).Value = -195530  # N77: was -165981
# Row 88
$ws.Cells.Item(88, 8).Value = 2789.625  # H88: was 3517.6667
$ws.Cells.Item(88, 9).Value = 3582.4  # I88: was 4326.5
$ws.Cells.Item(88, 10).Value = 1468.3334  # J88: was 1900
$ws.Cells.Item(88, 11).Value = 3582.4  # K88: was 4326.5
$ws.Cells.Item(88, 12).Value = 1468.3334  # L88: was 1900
$ws.Cells.Item(88, 13).Value = -3176.4  # M88: was -3920.5
$ws.Cells.Item(88, 14).Value = -2280.3334  # N88: was -2712
# Row 91
$ws.Cells.Item(91, 8).Value = 2789.625  # H91: was 3517.6667
$ws.Cells.Item(91, 9).Value = 3582.4  # I91: was 4326.5
$ws.Cells.Item(91, 10).Value = 1468.3334  # J91: was 1900
$ws.Cells.Item(91, 11).Value = 3582.4  # K91: was 4326.5
$ws.Cells.Item(91, 12).Value = 1468.3334  # L91: was 1900
$ws.Cells.Item(91, 13).Value = -2178.4  # M91: was -2922.5
$ws.Cells.Item(91, 14).Value = -4276.3334  # N91: was -4708
# Row 96
$ws.Cells.Item(96, 8).Value = 87500  # H96: was 0
$ws.Cells.Item(96, 10).Value = 87500  # J96: was 0
$ws.Cells.Item(96, 12).Value = 87500  # L96: was 0
$ws.Cells.Item(96, 14).Value = -92992  # N96: was None

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 1399.0555  # H94: was 1469.4706
$ws.Cells.Item(94, 9).Value = 1209.3334  # I94: was 1281.2858
$ws.Cells.Item(94, 11).Value = 1209.3334  # K94: was 1281.2858
$ws.Cells.Item(94, 13).Value = -758.3334  # M94: was -830.2858000000001
# Row 107
$ws.Cells.Item(107, 8).Value = 1145.6129  # H107: was 1229.6552
$ws.Cells.Item(107, 9).Value = 1037.6957  # I107: was 1116.409
$ws.Cells.Item(107, 10).Value = 1455.875  # J107: was 1585.5714
$ws.Cells.Item(107, 11).Value = 1037.6957  # K107: was 1116.409
$ws.Cells.Item(107, 12).Value = 1455.875  # L107: was 1585.5714
$ws.Cells.Item(107, 13).Value = 882.3043  # M107: was 803.5909999999999
$ws.Cells.Item(107, 14).Value = -5295.875  # N107: was -5425.5714
# Row 134
$ws.Cells.Item(134, 8).Value = 2035.875  # H134: was 2045.9
$ws.Cells.Item(134, 9).Value = 1882.8788  # I134: was 1910.5
$ws.Cells.Item(134, 10).Value = 2757.1428  # J134: was 2587.5
$ws.Cells.Item(134, 11).Value = 5648.636399999999  # K134: was 5731.5
$ws.Cells.Item(134, 12).Value = 8271.428400000001  # L134: was 7762.5
$ws.Cells.Item(134, 13).Value = -3113.636399999999  # M134: was -3196.5
$ws.Cells.Item(134, 14).Value = -13341.4284  # N134: was -12832.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 2666.3333  # H16: was 1437.6666
$ws.Cells.Item(16, 9).Value = 2666.3333  # I16: was 800
$ws.Cells.Item(16, 10).Value = 0  # J16: was 1756.5
$ws.Cells.Item(16, 11).Value = 2666.3333  # K16: was 800
$ws.Cells.Item(16, 12).Value = 0  # L16: was 1756.5
$ws.Cells.Item(16, 13).Value = -2379.3333  # M16: was -513
$ws.Cells.Item(16, 14).ClearContents()  # N16: was -2330.5
# Row 99
$ws.Cells.Item(99, 8).Value = 3093.75  # H99: was 3416.6667
$ws.Cells.Item(99, 9).Value = 2700  # I99: was 2875
$ws.Cells.Item(99, 10).Value = 3750  # J99: was 4500
$ws.Cells.Item(99, 11).Value = 2700  # K99: was 2875
$ws.Cells.Item(99, 12).Value = 3750  # L99: was 4500
$ws.Cells.Item(99, 13).Value = -1202  # M99: was -1377
$ws.Cells.Item(99, 14).Value = -6746  # N99: was -7496
# Row 105
$ws.Cells.Item(105, 8).Value = 1907.1428  # H105: was 1893.3334
$ws.Cells.Item(105, 9).Value = 1680  # I105: was 1602.5
$ws.Cells.Item(105, 11).Value = 1680  # K105: was 1602.5
$ws.Cells.Item(105, 13).Value = 67  # M105: was 144.5
# Row 107
$ws.Cells.Item(107, 8).Value = 356.2143  # H107: was 367.07693
$ws.Cells.Item(107, 9).Value = 360.53845  # I107: was 372.66666
$ws.Cells.Item(107, 11).Value = 360.53845  # K107: was 372.66666
$ws.Cells.Item(107, 13).Value = 1559.46155  # M107: was 1547.33334
# Row 113
$ws.Cells.Item(113, 8).Value = 2666.3333  # H113: was 1437.6666
$ws.Cells.Item(113, 9).Value = 2666.3333  # I113: was 800
$ws.Cells.Item(113, 10).Value = 0  # J113: was 1756.5
$ws.Cells.Item(113, 11).Value = 2666.3333  # K113: was 800
$ws.Cells.Item(113, 12).Value = 0  # L113: was 1756.5
$ws.Cells.Item(113, 13).Value = -496.3332999999998  # M113: was 1370
$ws.Cells.Item(113, 14).ClearContents()  # N113: was -6096.5
# Row 126
$ws.Cells.Item(126, 8).Value = 3093.75  # H126: was 3416.6667
$ws.Cells.Item(126, 9).Value = 2700  # I126: was 2875
$ws.Cells.Item(126, 10).Value = 3750  # J126: was 4500
$ws.Cells.Item(126, 11).Value = 8100  # K126: was 8625
$ws.Cells.Item(126, 12).Value = 11250  # L126: was 13500
$ws.Cells.Item(126, 13).Value = -5630  # M126: was -6155
$ws.Cells.Item(126, 14).Value = -16190  # N126: was -18440

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 93
$ws.Cells.Item(93, 8).Value = 6250  # H93: was 6500
$ws.Cells.Item(93, 10).Value = 6666.6665  # J93: was 8000
$ws.Cells.Item(93, 12).Value = 19999.9995  # L93: was 24000
$ws.Cells.Item(93, 14).Value = -23743.9995  # N93: was -27744
# Row 137
$ws.Cells.Item(137, 8).Value = 3133.1  # H137: was 3500.875
$ws.Cells.Item(137, 9).Value = 2070.4  # I137: was 2507
$ws.Cells.Item(137, 10).Value = 4195.8  # J137: was 4494.75
$ws.Cells.Item(137, 11).Value = 6211.200000000001  # K137: was 7521
$ws.Cells.Item(137, 12).Value = 12587.4  # L137: was 13484.25
$ws.Cells.Item(137, 13).Value = -1111.200000000001  # M137: was -2421
$ws.Cells.Item(137, 14).Value = -22787.4  # N137: was -23684.25
# Row 138
$ws.Cells.Item(138, 8).Value = 14714156  # H138: was 16675577
$ws.Cells.Item(138, 10).Value = 7277.577  # J138: was 7964.409
$ws.Cells.Item(138, 12).Value = 21832.731  # L138: was 23893.227
$ws.Cells.Item(138, 14).Value = -32112.731  # N138: was -34173.227

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Cells.Item(107, 8).Value = 1303.75  # H107: was 1407.5
$ws.Cells.Item(107, 9).Value = 1273.3334  # I107: was 1420
$ws.Cells.Item(107, 11).Value = 1273.3334  # K107: was 1420
$ws.Cells.Item(107, 13).Value = 646.6666  # M107: was 500
# Row 122
$ws.Cells.Item(122, 8).Value = 993.9231  # H122: was 1013.26086
$ws.Cells.Item(122, 9).Value = 951.0952  # I122: was 968.6667
$ws.Cells.Item(122, 11).Value = 2853.2856  # K122: was 2906.0001
$ws.Cells.Item(122, 13).Value = -403.2856000000002  # M122: was -456.0001000000002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 9998.083000000001  # H7: was 12916
$ws.Cells.Item(7, 9).Value = 11525.777  # I7: was 14463.429
$ws.Cells.Item(7, 10).Value = 5415  # J7: was 7500
$ws.Cells.Item(7, 11).Value = 11525.777  # K7: was 14463.429
$ws.Cells.Item(7, 12).Value = 5415  # L7: was 7500
$ws.Cells.Item(7, 13).Value = -11413.777  # M7: was -14351.429
$ws.Cells.Item(7, 14).Value = -5639  # N7: was -7724
# Row 40
$ws.Cells.Item(40, 8).Value = 6184.875  # H40: was 7501.25
$ws.Cells.Item(40, 9).Value = 5396.8  # I40: was 6000
$ws.Cells.Item(40, 10).Value = 7498.3335  # J40: was 9002.5
$ws.Cells.Item(40, 11).Value = 5396.8  # K40: was 6000
$ws.Cells.Item(40, 12).Value = 7498.3335  # L40: was 9002.5
$ws.Cells.Item(40, 13).Value = -5260.8  # M40: was -5864
$ws.Cells.Item(40, 14).Value = -7770.3335  # N40: was -9274.5
# Row 55
$ws.Cells.Item(55, 8).Value = 329.6316  # H55: was 328.05264
$ws.Cells.Item(55, 9).Value = 300.0909  # I55: was 297.36365
$ws.Cells.Item(55, 11).Value = 300.0909  # K55: was 297.36365
$ws.Cells.Item(55, 13).Value = -127.0909  # M55: was -124.36365
# Row 61
$ws.Cells.Item(61, 8).Value = 2921.348  # H61: was 3124
$ws.Cells.Item(61, 9).Value = 3210.647  # I61: was 3532.9333
$ws.Cells.Item(61, 11).Value = 3210.647  # K61: was 3532.9333
$ws.Cells.Item(61, 13).Value = -3008.647  # M61: was -3330.9333
# Row 93
$ws.Cells.Item(93, 8).Value = 2613.111  # H93: was 2846
$ws.Cells.Item(93, 9).Value = 2567.9375  # I93: was 2827.6428
$ws.Cells.Item(93, 11).Value = 2567.9375  # K93: was 2827.6428
$ws.Cells.Item(93, 13).Value = -1319.9375  # M93: was -1579.6428
# Row 113
$ws.Cells.Item(113, 8).Value = 2921.348  # H113: was 3124
$ws.Cells.Item(113, 9).Value = 3210.647  # I113: was 3532.9333
$ws.Cells.Item(113, 11).Value = 3210.647  # K113: was 3532.9333
$ws.Cells.Item(113, 13).Value = -1040.647  # M113: was -1362.9333
# Row 122
$ws.Cells.Item(122, 8).Value = 5722.222  # H122: was 5450
$ws.Cells.Item(122, 9).Value = 4785.7144  # I122: was 4562.5
$ws.Cells.Item(122, 11).Value = 14357.1432  # K122: was 13687.5
$ws.Cells.Item(122, 13).Value = -11907.1432  # M122: was -11237.5
# Row 126
$ws.Cells.Item(126, 8).Value = 9998.083000000001  # H126: was 12916
$ws.Cells.Item(126, 9).Value = 11525.777  # I126: was 14463.429
$ws.Cells.Item(126, 10).Value = 5415  # J126: was 7500
$ws.Cells.Item(126, 11).Value = 34577.331  # K126: was 43390.287
$ws.Cells.Item(126, 12).Value = 16245  # L126: was 22500
$ws.Cells.Item(126, 13).Value = -32107.331  # M126: was -40920.287
$ws.Cells.Item(126, 14).Value = -21185  # N126: was -27440
# Row 136
$ws.Cells.Item(136, 8).Value = 3780.4814  # H136: was 4213.6523
$ws.Cells.Item(136, 9).Value = 3048.818  # I136: was 3439.7222
$ws.Cells.Item(136, 11).Value = 9146.454000000002  # K136: was 10319.1666
$ws.Cells.Item(136, 13).Value = -6596.454000000002  # M136: was -7769.1666

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 28
$ws.Cells.Item(28, 8).Value = 17925  # H28: was 16698.334
$ws.Cells.Item(28, 10).Value = 17925  # J28: was 16698.334
$ws.Cells.Item(28, 12).Value = 17925  # L28: was 16698.334
$ws.Cells.Item(28, 14).Value = -18621  # N28: was -17394.334
# Row 122
$ws.Cells.Item(122, 8).Value = 2656  # H122: was 2175.5945
$ws.Cells.Item(122, 9).Value = 2456  # I122: was 2028.1515
$ws.Cells.Item(122, 10).Value = 4322.6665  # J122: was 3392
$ws.Cells.Item(122, 11).Value = 7368  # K122: was 6084.4545
$ws.Cells.Item(122, 12).Value = 12967.9995  # L122: was 10176
$ws.Cells.Item(122, 13).Value = -4918  # M122: was -3634.4545
$ws.Cells.Item(122, 14).Value = -17867.9995  # N122: was -15076
# Row 132
$ws.Cells.Item(132, 8).Value = 2624.8958  # H132: was 2961.8298
$ws.Cells.Item(132, 9).Value = 2336.0264  # I132: was 2651.1538
$ws.Cells.Item(132, 10).Value = 3722.6  # J132: was 4476.375
$ws.Cells.Item(132, 11).Value = 7008.0792  # K132: was 7953.4614
$ws.Cells.Item(132, 12).Value = 11167.8  # L132: was 13429.125
$ws.Cells.Item(132, 13).Value = -4478.0792  # M132: was -5423.4614
$ws.Cells.Item(132, 14).Value = -16227.8  # N132: was -18489.125
